$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for records 32-35 got cyclically rotated: the content that
# used to live in row 33 now lives in row 32, row 34's content moved to row
# 33, row 35's content moved to row 34, and the original row 32 content
# (which carried some extra columns: J, L, M, N, AF) moved down to row 35.

# --- Row 32 (becomes former row 33's record) ---
$ws.Range("A32").Value = 111666918
$ws.Range("B32").Value = 77550
$ws.Range("D32").Value = "NT"
$ws.Range("E32").Value = 185
$ws.Range("F32").Value = "Violettgrå tagellav"
$ws.Range("G32").Value = "Bryoria nadvornikiana"
$ws.Range("H32").Value = "(Gyeln.) Brodo & D.Hawksw."
$ws.Range("P32").Value = "Svartflärksbäcken (Svartflärksbäcken), Mpd"
$ws.Range("Q32").Value = 627982
$ws.Range("R32").Value = 6943734
# Row 32 no longer carries these optional columns (they belonged to the
# original row-32 record, which has now moved down to row 35).
$ws.Range("J32").ClearContents()
$ws.Range("L32").ClearContents()
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()
$ws.Range("AF32").ClearContents()

# --- Row 33 (becomes former row 34's record) ---
$ws.Range("A33").Value = 111668313
$ws.Range("B33").Value = 56543
$ws.Range("D33").Value = "NT"
$ws.Range("E33").Value = 103021
$ws.Range("F33").Value = "Talltita"
$ws.Range("G33").Value = "Poecile montanus"
$ws.Range("H33").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("P33").Value = "Svartflärken (Svartflärken), Mpd"
$ws.Range("Q33").Value = 627870
$ws.Range("R33").Value = 6944135

# --- Row 34 (becomes former row 35's record) ---
$ws.Range("A34").Value = 111667877
$ws.Range("B34").Value = 56543
$ws.Range("D34").Value = "NT"
$ws.Range("E34").Value = 103021
$ws.Range("F34").Value = "Talltita"
$ws.Range("G34").Value = "Poecile montanus"
$ws.Range("H34").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("P34").Value = "Svartflärken (Svartflärken), Mpd"
$ws.Range("Q34").Value = 628073
$ws.Range("R34").Value = 6944660

# --- Row 35 (becomes former row 32's record) ---
$ws.Range("A35").Value = 111668109
$ws.Range("B35").Value = 5135
$ws.Range("D35").Value = "LC"
$ws.Range("E35").Value = 105930
$ws.Range("F35").Value = "Vågbandad barkbock"
$ws.Range("G35").Value = "Semanotus undatus"
$ws.Range("H35").Value = "(Linnaeus, 1758)"
$ws.Range("P35").Value = "Svartflärken (Svartflärken), Mpd"
$ws.Range("Q35").Value = 628016
$ws.Range("R35").Value = 6944481
# Row 35 now carries the optional columns that used to belong to row 32.
$ws.Range("M35").Value = "äldre gnagspår"
